$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Agosto de 2020 a las 16:32"

# --- Update per-country statistics (columns B:H) ---
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 5570145
$ws.Range("C4").Value = 3513
$ws.Range("D4").Value = 2922964
$ws.Range("E4").Value = 2474033
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = 173148

# Row 6 - India
$ws.Range("B6").Value = 2667973
$ws.Range("C6").Value = 20657
$ws.Range("D6").Value = 1939454
$ws.Range("E6").Value = 677239
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 235
$ws.Range("H6").Value = 51280

# Row 15 - Reino Unido
$ws.Range("B15").Value = 319197
$ws.Range("C15").Value = 713
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 41369

# Row 22 - Alemania
$ws.Range("B22").Value = 225696
$ws.Range("C22").Value = 699
$ws.Range("D22").Value = 202900
$ws.Range("E22").Value = 13504
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 9292

# Row 28 - Catar
$ws.Range("B28").Value = 115368
$ws.Range("C28").Value = 288
$ws.Range("D28").Value = 112088
$ws.Range("E28").Value = 3087
$ws.Range("F28").Value = 0

# Row 43 - Bielorrusia
$ws.Range("B43").Value = 69589
$ws.Range("C43").Value = 73
$ws.Range("D43").Value = 67149
$ws.Range("E43").Value = 1827
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 3
$ws.Range("H43").Value = 613

# Row 62 - Azerbaiyan
$ws.Range("B62").Value = 34343
$ws.Range("C62").Value = 124
$ws.Range("D62").Value = 32042
$ws.Range("E62").Value = 1793
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 508

# Row 67 - Serbia
$ws.Range("B67").Value = 29782
$ws.Range("C67").Value = 100
$ws.Range("D67").Value = 27208
$ws.Range("E67").Value = 1897
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 3
$ws.Range("H67").Value = 677

# Row 77 - Estado de Palestina
$ws.Range("B77").Value = 16844
$ws.Range("C77").Value = 310
$ws.Range("D77").Value = 9838
$ws.Range("E77").Value = 6896
$ws.Range("F77").Value = 0

# Row 86 - Noruega
$ws.Range("B86").Value = 10035
$ws.Range("C86").Value = 30
$ws.Range("D86").Value = 8857
$ws.Range("E86").Value = 917
$ws.Range("F86").Value = 0

# Row 88 - Consejo Danes para los Refugiados
$ws.Range("B88").Value = 9706
$ws.Range("C88").Value = 30
$ws.Range("D88").Value = 8705
$ws.Range("E88").Value = 758
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 3
$ws.Range("H88").Value = 243

# --- Montserrat / Islas Malvinas swap places in the list (rows 213 & 214) ---
# Row 213 was Montserrat (13,0,12,0,0,0,1) -> becomes Islas Malvinas (13,0,13,0,0,0,0)
# Row 214 was Islas Malvinas (13,0,13,0,0,0,0) -> becomes Montserrat (13,0,12,0,0,0,1)
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 13
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1
